$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-9) were reordered. Columns A,B,C,E,F,G,H,I,J,Q,T are
# identical across all rows, so only D (Fecha), K (Variedad), L (Calidad),
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado), R (Origen) and S (Precio $/Kg) need to be rewritten to match
# the new row order.

$data = @{
    2  = @{ D = 45043; K = "Fuyu";     L = "Primera"; M = 300; N = 25000; O = 26000; P = 25500; R = "Región de O'Higgins"; S = 1417 }
    3  = @{ D = 44305; K = "Mankaki";  L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    4  = @{ D = 44342; K = "Mankaki";  L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    5  = @{ D = 44313; K = "Mankaki";  L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 }
    6  = @{ D = 44301; K = "Hachiya";  L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 }
    7  = @{ D = 44355; K = "Mankaki";  L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana"; S = 1139 }
    8  = @{ D = 44699; K = "Mankaki";  L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins"; S = 1639 }
    9  = @{ D = 45071; K = "Fuyu";     L = "Segunda"; M = 110; N = 23000; O = 24000; P = 23455; R = "Región Metropolitana"; S = 1303 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 11).Value = $vals.K   # K - Variedad
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $vals.R   # R - Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S - Precio $/Kg
}
